$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.044.32"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "3.862.61"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "696.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.47%  "

$ws.Range("D7").Value = "3.861.62"
$ws.Range("E7").Value = "  +1.44%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.05%  "

$ws.Range("E12").Value = "  -0.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("D15").Value = "4.517.76"
$ws.Range("E15").Value = "  +1.54%  "

$ws.Range("D16").Value = "3.869.88"
$ws.Range("E16").Value = "  +1.81%  "

$ws.Range("D17").Value = "71.095.89"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "497.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.721"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.60%  "

$ws.Range("E24").Value = "  +3.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.12%  "

$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.86%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("E34").Value = "  +2.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("D36").Value = "3.820.86"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.104"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.68%  "

$ws.Range("E40").Value = "  +8.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("E45").Value = "  -7.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.36"
$ws.Range("D47").Style = "Normal"

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "417.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.62%  "

$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.302"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("E50").Value = "  -1.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.33%  "
